# Update cryptocurrency price/volume data (and a couple of row re-orderings)
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.333.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.35%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.842.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.43%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9978'
$ws.Range('D4').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.07%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6291'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07455'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.48%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2896'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.52%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.95%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07725'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.21%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.842.47'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.42%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.998'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.73%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6783'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.38%  '

# Row 15
$ws.Range('E15').Value = '  -4.63%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.71%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.135'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.80%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.357.20'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.56%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.18%  '

# Row 20
$ws.Range('E20').Value = '  -0.26%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9988'
$ws.Range('D21').Style = 'Normal'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.417'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.52%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9991'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.27%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.78%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1375'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.84%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.414'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.25%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.94%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06455'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +15.09%  '

# Row 29
$ws.Range('E29').Value = '  +0.49%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.474'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.82%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.083'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.23%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.049'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.11%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.823'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.29%  '

# Row 34
$ws.Range('E34').Value = '  -1.95%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6954'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.00%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.581'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.35%  '

# Row 37
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.256.74'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.09%  '

# Row 38
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.829'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.90%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.519'
$ws.Range('D40').Style = 'Normal'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9095'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.27%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9981'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.34%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.004.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -13.85%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.79%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.35'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.37%  '

# Row 46
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1169'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '

# Row 47
$ws.Range('E47').Value = '  -2.36%  '

# Row 48
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000115'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.69%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.001'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.12%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3945'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.01%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.675'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
